$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy style from existing header cell (H1) to new headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

$iValues = @(9,6,6,7,1,7,5,11,8,9,7,4,8,6,7,5,9,9,6,8,4,7,6,4,8,7,7,7,7,6,5,4,6,3,6)
$jValues = @(9,6,6,8,1,7,5,11,8,9,7,5,8,6,7,6,9,9,6,8,5,7,7,5,8,7,8,7,7,6,5,5,6,3,6)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
